# Update serial-number values in the flapdoodle Items sheet.
# SN990005 -> SN990008, SN990006 -> SN990009, SN990007 -> SN990010
# These values appear in a small circular "self-subcomponent" chain:
#   A8/C10 = SN990005->SN990008
#   A9/C8  = SN990006->SN990009
#   A10/C9 = SN990007->SN990010

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date system stays 1900-based (unchanged), kept explicit for clarity.
$wb.Date1904 = $false

$ws.Range("A8").Value = "SN990008"
$ws.Range("C10").Value = "SN990008"

$ws.Range("A9").Value = "SN990009"
$ws.Range("C8").Value = "SN990009"

$ws.Range("A10").Value = "SN990010"
$ws.Range("C9").Value = "SN990010"

# Move the active selection, matching the author's final cursor position.
$ws.Range("D9").Select()
